# Fruta / hortaliza, semanal
# Insert a new weekly record at row 6 (pushing all subsequent records down
# by one row) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6:48 down to 7:49, duplicating row 6's formatting (incl. the
# date number-format on column D) into the freshly inserted row.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record.
$ws.Cells.Item(6, 1).Value  = 11
$ws.Cells.Item(6, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(6, 3).Value  = "Bíobío"
$ws.Cells.Item(6, 4).Value  = 44490
$ws.Cells.Item(6, 5).Value  = 8
$ws.Cells.Item(6, 6).Value  = 100112021
$ws.Cells.Item(6, 7).Value  = "Ají"
$ws.Cells.Item(6, 8).Value  = "Americana (o)"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 50
$ws.Cells.Item(6, 11).Value = 62000
$ws.Cells.Item(6, 12).Value = 64000
$ws.Cells.Item(6, 13).Value = 63200
$ws.Cells.Item(6, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 2528
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
